$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.041.51"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.565.57"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.50"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.492"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.13"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0597"
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0858"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.562.60"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.77"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.023.45"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.87"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0706"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "216.02"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.42"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("E21").Value = "  +2.27%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "154.11"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.61"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.06"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("E27").Value = "  +1.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0475"
$ws.Range("E29").Value = "  +1.74%  "
$ws.Range("E30").Value = "  +3.88%  "
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  +4.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.427.49"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("E34").Value = "  +12.37%  "
$ws.Range("E35").Value = "  +1.95%  "
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.535"
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.01"
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.33"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.82"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.700.64"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.56"
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("E48").Value = "  +1.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0520"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0961"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("E51").Value = "  +0.25%  "
